$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.379247069358826
$ws.Range("B1").Value = 2.696710348129272
$ws.Range("C1").Value = 3.296211957931519
$ws.Range("D1").Value = 3.290188312530518
$ws.Range("E1").Value = 2.014054298400879
